# Refresh crypto price/volume/coin data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''60.644.87'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.47%  '

$ws.Range('D3').Value = '''2.399.46'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.85%  '

$ws.Range('E4').Value = '  +0.46%  '

$ws.Range('D5').Value = '''563.76'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.33%  '

$ws.Range('D6').Value = '''141.30'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.89%  '

$ws.Range('E8').Value = '  +1.29%  '

$ws.Range('D9').Value = '''2.406.21'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('E10').Value = '  -0.25%  '

$ws.Range('E11').Value = '  -0.94%  '

$ws.Range('D12').Value = '''5.17'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.51%  '

$ws.Range('E13').Value = '  -0.10%  '

$ws.Range('D14').Value = '''26.10'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.46%  '

$ws.Range('D15').Value = '''0.0000168'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.99%  '

$ws.Range('D16').Value = '''2.822.80'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.15%  '

$ws.Range('D17').Value = '''60.580.54'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.35%  '

$ws.Range('D18').Value = '''2.407.15'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.44%  '

$ws.Range('D19').Value = '''8.03'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +8.26%  '

$ws.Range('D20').Value = '''10.63'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.01%  '

$ws.Range('D21').Value = '''324.29'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.40%  '

$ws.Range('E22').Value = '  +0.62%  '

$ws.Range('E23').Value = '  -0.18%  '

$ws.Range('E24').Value = '  +0.01%  '

$ws.Range('D25').Value = '''1.83'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.66%  '

$ws.Range('D26').Value = '''65.04'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.06%  '

$ws.Range('D27').Value = '''569.60'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.28%  '

$ws.Range('D28').Value = '''8.06'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.38%  '

$ws.Range('D29').Value = '''2.513.61'
$ws.Range('D29').ClearFormats()

$ws.Range('D30').Value = '''0.0₃0937'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.42%  '

$ws.Range('E31').Value = '  +1.95%  '

$ws.Range('E32').Value = '  -1.48%  '

$ws.Range('E33').Value = '  -1.42%  '

$ws.Range('E34').Value = '  -1.80%  '

$ws.Range('E35').Value = '  -0.51%  '

$ws.Range('E36').Value = '  +3.90%  '

$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = '''153.99'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.22%  '

$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').Value = '''0.371'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.53%  '

$ws.Range('E39').Value = '  -1.53%  '

$ws.Range('D40').Value = '''18.30'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.04%  '

$ws.Range('D41').Value = '''5.14'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.33%  '

$ws.Range('E42').Value = '  -0.08%  '

$ws.Range('D43').Value = '''2.51'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +7.52%  '

$ws.Range('E44').Value = '  +0.34%  '

$ws.Range('D45').Value = '''41.66'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.28%  '

$ws.Range('E46').Value = '  +3.43%  '

$ws.Range('D47').Value = '''141.64'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.23%  '

$ws.Range('D48').Value = '''3.56'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.80%  '

$ws.Range('E49').Value = '  +0.16%  '

$ws.Range('E50').Value = '  +0.26%  '

$ws.Range('D51').Value = '''19.34'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.23%  '
